$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Copy formatting from column K (last existing data column) into the new
# columns L and M for rows 7-102, so the new cells inherit the right number formats/styles
# (date format for row 7/38/80, numeric format elsewhere) before we fill in values.
$ws.Range("K7:K102").Copy()
$ws.Range("L7:M7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 2: Write the updated quarterly figures.
$ws.Range("D7").Value = 43496
$ws.Range("E7").Value = 43404
$ws.Range("F7").Value = 43312
$ws.Range("G7").Value = 43220
$ws.Range("H7").Value = 43131
$ws.Range("I7").Value = 43039
$ws.Range("J7").Value = 42947
$ws.Range("K7").Value = 42855
$ws.Range("L7").Value = 42766
$ws.Range("M7").Value = 42674
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 0
$ws.Range("L9").Value = "NA"
$ws.Range("M9").Value = "NA"
$ws.Range("L10").Value = "NA"
$ws.Range("M10").Value = "NA"
$ws.Range("L12").Value = "NA"
$ws.Range("M12").Value = "NA"
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("L21").Value = "NA"
$ws.Range("M21").Value = "NA"
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = 0
$ws.Range("D38").Value = 43496
$ws.Range("E38").Value = 43404
$ws.Range("F38").Value = 43312
$ws.Range("G38").Value = 43220
$ws.Range("H38").Value = 43131
$ws.Range("I38").Value = 43039
$ws.Range("J38").Value = 42947
$ws.Range("K38").Value = 42855
$ws.Range("L38").Value = 42766
$ws.Range("M38").Value = 42674
$ws.Range("E41").Value = 0
$ws.Range("F41").Value = 0
$ws.Range("L41").Value = 200
$ws.Range("M41").Value = 200
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = 0
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "NA"
$ws.Range("H43").Value = "NA"
$ws.Range("I43").Value = "NA"
$ws.Range("J43").Value = "NA"
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("L46").Value = 200
$ws.Range("M46").Value = 200
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 0
$ws.Range("E54").Value = 0
$ws.Range("F54").Value = 0
$ws.Range("L54").Value = 200
$ws.Range("M54").Value = 200
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = 0
$ws.Range("K72").Value = -2200
$ws.Range("L72").Value = -2200
$ws.Range("M72").Value = -2100
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = 0
$ws.Range("E76").Value = 0
$ws.Range("F76").Value = 0
$ws.Range("L76").Value = 200
$ws.Range("M76").Value = 200
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 0
$ws.Range("D80").Value = 43496
$ws.Range("E80").Value = 43404
$ws.Range("F80").Value = 43312
$ws.Range("G80").Value = 43220
$ws.Range("H80").Value = 43131
$ws.Range("I80").Value = 43039
$ws.Range("J80").Value = 42947
$ws.Range("K80").Value = 42855
$ws.Range("L80").Value = 42766
$ws.Range("M80").Value = 42674
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("F96").Value = -200
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 0
$ws.Range("F100").Value = -200
$ws.Range("L100").Value = "NA"
$ws.Range("M100").Value = "NA"
$ws.Range("L101").Value = 0
$ws.Range("M101").Value = 0
$ws.Range("D102").Value = 0
$ws.Range("F102").Value = -200
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 0
